# Update "Price" (column D) and "Volume(1h)" (column E) for the cryptos
# list as refreshed by the scheduled GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NumberFormat is forced to Text ("@") before assigning numeric-looking
# strings (e.g. "1.009", "0.06850") so Excel keeps them as literal text
# instead of silently re-parsing them as numbers and dropping trailing
# zeros / precision. The style is reset back to "Normal" afterwards so
# the cell's formatting stays the same as before the edit.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.513.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.92%  "

$ws.Range("E4").Value = "  -0.75%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5084"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3905"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.103"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.230"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.879.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.266"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.009"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001104"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06734"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("E19").Value = "  +0.94%  "

$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.530.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.236"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.082.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.390"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "125.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.1046"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.040"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.777"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.612"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.02457"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06546"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2160"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.816"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.059"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.08%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.255"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.192"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6399"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("E43").Value = "  -0.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6017"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.689"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.006"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("E48").Value = "  +1.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.135"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -11.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06850"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.12%  "
